# Apply the "updated 4.0 files and mdl" commit to the Maximum Capacity Factor
# workbook:
#   - "About" sheet: bump the last-updated date in C1 (45320 -> 45392,
#     i.e. 1/29/2024 -> 4/10/2024)
#   - "MCF" sheet: raise every non-variable/non-zero capacity factor to 1
#     (the dependent formula cells B19:B25 recalc automatically since they
#     just reference the edited cells)
#   - leave the active-cell selection on the MCF sheet at B17, matching the
#     state the workbook was saved in

$wb = $excel.ActiveWorkbook

# --- "About" sheet: update the last-modified date stamp ---
$about = $wb.Worksheets.Item("About")
$about.Range("C1").Value = 45392

# --- "MCF" sheet: bump capacity factors that were 0.85 / 0.95 up to 1 ---
$mcf = $wb.Worksheets.Item("MCF")

$mcf.Range("B2").Value = 1    # hard coal
$mcf.Range("B3").Value = 1    # natural gas steam turbine
$mcf.Range("B4").Value = 1    # natural gas combined cycle
$mcf.Range("B6").Value = 1    # hydro
$mcf.Range("B10").Value = 1   # biomass
$mcf.Range("B11").Value = 1   # geothermal
$mcf.Range("B12").Value = 1   # petroleum
$mcf.Range("B13").Value = 1   # natural gas peaker
$mcf.Range("B14").Value = 1   # lignite
$mcf.Range("B16").Value = 1   # crude oil
$mcf.Range("B17").Value = 1   # heavy or residual fuel oil
$mcf.Range("B18").Value = 1   # municipal solid waste

# B19:B25 are formulas referencing the cells above (=B2, =B4, =B10, =B14,
# =B5, =B4, =B4) and recalculate automatically.

# Restore the sheet's active cell/selection to match the saved state.
$mcf.Activate()
$mcf.Range("B17").Select()
